$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 1 above the existing data, shifting everything down by one.
$ws.Rows.Item(1).Insert()

# Write the new header row.
$ws.Range("A1").Value = "Município"
$ws.Range("B1").Value = "Confirmados"
$ws.Range("C1").Value = "Óbitos"

# Style the header row: bold font, thin border all around, centered / top aligned.
$hdr = $ws.Range("A1:C1")
$hdr.Font.Bold = $true
$hdr.Borders.LineStyle = 1
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
